# IACET_Dispatcher folder added and in this folder read mail and check
# attachment workflow added -- add the new mail-folder / subject
# configuration entries to the "Settings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New Name/Value/Description rows describing the mail folders & subjects
# used by the new IACET_Dispatcher workflow.
$ws.Range("A6").Value = "FolderName"
$ws.Range("B6").Value = "Inbox"

$ws.Range("A7").Value = "ProjectSubjectFolder"
$ws.Range("B7").Value = "IACET"

$ws.Range("A8").Value = "OtherSubjectFolder"
$ws.Range("B8").Value = "Other mails"

$ws.Range("A9").Value = "ProjectSubject"
$ws.Range("B9").Value = "IACET client Details"

# Make "Settings" the active sheet/tab with the same selection Excel left
# it in when the author saved the workbook.
$ws.Activate()
$ws.Range("A14").Select()
